$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 147; Date = "02-08-2021"; B = 13264; C = 35382; D = 214; E = 26761 },
    @{ Row = 148; Date = "03-08-2021"; B = 13708; C = 35426; D = 214; E = 26814 },
    @{ Row = 149; Date = "04-08-2021"; B = 12861; C = 35467; D = 214; E = 26640 },
    @{ Row = 150; Date = "05-08-2021"; B = 13564; C = 35570; D = 214; E = 26774 },
    @{ Row = 151; Date = "06-08-2021"; B = 14283; C = 35584; D = 214; E = 27169 },
    @{ Row = 152; Date = "09-08-2021"; B = 13707; C = 35604; D = 214; E = 27341 },
    @{ Row = 153; Date = "10-08-2021"; B = 13903; C = 35616; D = 214; E = 27306 },
    @{ Row = 154; Date = "11-08-2021"; B = 13642; C = 35319; D = 214; E = 27527 },
    @{ Row = 155; Date = "12-08-2021"; B = 14591; C = 35248; D = 214; E = 27780 },
    @{ Row = 156; Date = "13-08-2021"; B = 15328; C = 35778; D = 214; E = 28038 },
    @{ Row = 157; Date = "16-08-2021"; B = 14656; C = 35746; D = 214; E = 28184 },
    @{ Row = 158; Date = "17-08-2021"; B = 15204; C = 35776; D = 214; E = 28239 },
    @{ Row = 159; Date = "18-08-2021"; B = 15282; C = 35753; D = 214; E = 28280 },
    @{ Row = 160; Date = "19-08-2021"; B = 15399; C = 35775; D = 214; E = 28314 },
    @{ Row = 161; Date = "20-08-2021"; B = 16045; C = 36137; D = 214; E = 28064 },
    @{ Row = 162; Date = "23-08-2021"; B = 15353; C = 36184; D = 114; E = 27957 },
    @{ Row = 163; Date = "24-08-2021"; B = 15385; C = 36178; D = 114; E = 27982 },
    @{ Row = 164; Date = "25-08-2021"; B = 17016; C = 36147; D = 114; E = 28006 },
    @{ Row = 165; Date = "26-08-2021"; B = 15870; C = 36160; D = 114; E = 27664 },
    @{ Row = 166; Date = "27-08-2021"; B = 16045; C = 36142; D = 114; E = 27825 },
    @{ Row = 167; Date = "30-08-2021"; B = 16155; C = 36210; D = 114; E = 27850 },
    @{ Row = 168; Date = "31-08-2021"; B = 15608; C = 36238; D = 114; E = 28096 }
)

foreach ($row in $newRows) {
    $r = $row.Row

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Formula = "=""" + $row.Date + """"
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
}

$excel.CutCopyMode = 0

